$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2908.8462
$ws.Range("J64").Value = 3002.5
$ws.Range("L64").Value = 3002.5
$ws.Range("N64").Value = -3498.5

$ws.Range("H67").Value = 2908.8462
$ws.Range("J67").Value = 3002.5
$ws.Range("L67").Value = 3002.5
$ws.Range("N67").Value = -4718.5

$ws.Range("H131").Value = 3806
$ws.Range("J131").Value = 3806
$ws.Range("L131").Value = 11418
$ws.Range("N131").Value = -21498

$ws.Range("H136").Value = 60780
$ws.Range("J136").Value = 60780
$ws.Range("L136").Value = 60780
$ws.Range("N136").Value = -70980

$ws.Range("H137").Value = 2912.4092
$ws.Range("I137").Value = 3037.5
$ws.Range("J137").Value = 2578.8333
$ws.Range("K137").Value = 9112.5
$ws.Range("L137").Value = 7736.499899999999
$ws.Range("M137").Value = -6562.5
$ws.Range("N137").Value = -12836.4999

$ws.Range("H138").Value = 1956.1548
$ws.Range("I138").Value = 1488.2
$ws.Range("K138").Value = 4464.6
$ws.Range("M138").Value = 675.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 956606.2
$ws.Range("I32").Value = 1019763.25
$ws.Range("K32").Value = 1019763.25
$ws.Range("M32").Value = -1019476.25

$ws.Range("H45").Value = 3949.3333
$ws.Range("I45").Value = 4932
$ws.Range("K45").Value = 4932
$ws.Range("M45").Value = -4555

$ws.Range("H61").Value = 3244.4583
$ws.Range("I61").Value = 3153.182
$ws.Range("J61").Value = 3321.6924
$ws.Range("K61").Value = 3153.182
$ws.Range("L61").Value = 3321.6924
$ws.Range("M61").Value = -2941.182
$ws.Range("N61").Value = -3745.6924

$ws.Range("H74").Value = 905.23914
$ws.Range("I74").Value = 690.15625
$ws.Range("J74").Value = 1396.8572
$ws.Range("K74").Value = 690.15625
$ws.Range("L74").Value = 1396.8572
$ws.Range("M74").Value = 183.84375
$ws.Range("N74").Value = -3144.8572

$ws.Range("H77").Value = 905.23914
$ws.Range("I77").Value = 690.15625
$ws.Range("J77").Value = 1396.8572
$ws.Range("K77").Value = 3450.78125
$ws.Range("L77").Value = 6984.286
$ws.Range("M77").Value = 917.21875
$ws.Range("N77").Value = -15720.286

$ws.Range("H122").Value = 58184.61
$ws.Range("J122").Value = 2711.1428
$ws.Range("L122").Value = 8133.428400000001
$ws.Range("N122").Value = -13033.4284

$ws.Range("H123").Value = 32929
$ws.Range("J123").Value = 32929
$ws.Range("L123").Value = 32929
$ws.Range("N123").Value = -42729

$ws.Range("H132").Value = 2920.0747
$ws.Range("I132").Value = 2356.75
$ws.Range("J132").Value = 4343.2104
$ws.Range("K132").Value = 7070.25
$ws.Range("L132").Value = 13029.6312
$ws.Range("M132").Value = -4540.25
$ws.Range("N132").Value = -18089.6312

$ws.Range("H136").Value = 3244.4583
$ws.Range("I136").Value = 3153.182
$ws.Range("J136").Value = 3321.6924
$ws.Range("K136").Value = 9459.545999999998
$ws.Range("L136").Value = 9965.0772
$ws.Range("M136").Value = -6909.545999999998
$ws.Range("N136").Value = -15065.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 42000
$ws.Range("J111").Value = 42000
$ws.Range("L111").Value = 42000
$ws.Range("N111").Value = -50180

$ws.Range("H134").Value = 4296.25
$ws.Range("I134").Value = 3696.4614
$ws.Range("J134").Value = 5005.091
$ws.Range("K134").Value = 11089.3842
$ws.Range("L134").Value = 15015.273
$ws.Range("M134").Value = -8554.3842
$ws.Range("N134").Value = -20085.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3047.0688
$ws.Range("I31").Value = 951.8431399999999
$ws.Range("K31").Value = 951.8431399999999
$ws.Range("M31").Value = -656.8431399999999

$ws.Range("H34").Value = 3047.0688
$ws.Range("I34").Value = 951.8431399999999
$ws.Range("K34").Value = 951.8431399999999
$ws.Range("M34").Value = -749.8431399999999

$ws.Range("H132").Value = 2050.4546
$ws.Range("I132").Value = 1910.7084
$ws.Range("K132").Value = 5732.1252
$ws.Range("M132").Value = -3202.1252

$ws.Range("H134").Value = 4113.8945
$ws.Range("I134").Value = 4904.7856
$ws.Range("J134").Value = 1899.4
$ws.Range("K134").Value = 14714.3568
$ws.Range("L134").Value = 5698.200000000001
$ws.Range("M134").Value = -12179.3568
$ws.Range("N134").Value = -10768.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1910.2727
$ws.Range("I5").Value = 775.2857
$ws.Range("J5").Value = 2439.9333
$ws.Range("K5").Value = 2325.8571
$ws.Range("L5").Value = 7319.7999
$ws.Range("M5").Value = -2213.8571
$ws.Range("N5").Value = -7543.7999

$ws.Range("H92").Value = 1000
$ws.Range("I92").Value = 1000
$ws.Range("K92").Value = 3000
$ws.Range("M92").Value = -1752

$ws.Range("H135").Value = 1910.2727
$ws.Range("I135").Value = 775.2857
$ws.Range("J135").Value = 2439.9333
$ws.Range("K135").Value = 6977.571300000001
$ws.Range("L135").Value = 21959.3997
$ws.Range("M135").Value = -4442.571300000001
$ws.Range("N135").Value = -27029.3997

$ws.Range("H139").Value = 2555.568
$ws.Range("I139").Value = 1378.1364
$ws.Range("J139").Value = 3733
$ws.Range("K139").Value = 4134.4092
$ws.Range("L139").Value = 11199
$ws.Range("M139").Value = 1005.5908
$ws.Range("N139").Value = -21479

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2144.375
$ws.Range("I122").Value = 2196.7144
$ws.Range("J122").Value = 1778
$ws.Range("K122").Value = 6590.1432
$ws.Range("L122").Value = 5334
$ws.Range("M122").Value = -4140.1432
$ws.Range("N122").Value = -10234

$ws.Range("H123").Value = 10201.889
$ws.Range("J123").Value = 10201.889
$ws.Range("L123").Value = 10201.889
$ws.Range("N123").Value = -15101.889

$ws.Range("H132").Value = 4670.3335
$ws.Range("I132").Value = 4415.6665
$ws.Range("J132").Value = 4925
$ws.Range("K132").Value = 13246.9995
$ws.Range("L132").Value = 14775
$ws.Range("M132").Value = -10716.9995
$ws.Range("N132").Value = -19835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 35506.5
$ws.Range("I24").Value = 1006
$ws.Range("J24").Value = 70007
$ws.Range("K24").Value = 1006
$ws.Range("L24").Value = 70007
$ws.Range("M24").Value = -663
$ws.Range("N24").Value = -70693

$ws.Range("H122").Value = 3663.2307
$ws.Range("I122").Value = 3332.182
$ws.Range("J122").Value = 3906
$ws.Range("K122").Value = 9996.545999999998
$ws.Range("L122").Value = 11718
$ws.Range("M122").Value = -7546.545999999998
$ws.Range("N122").Value = -16618

$ws.Range("H136").Value = 4763651.5
$ws.Range("I136").Value = 1960.16
$ws.Range("J136").Value = 16667880
$ws.Range("K136").Value = 5880.48
$ws.Range("L136").Value = 50003640
$ws.Range("M136").Value = -3330.48
$ws.Range("N136").Value = -50008740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2482.8667
$ws.Range("I122").Value = 2504.2173
$ws.Range("K122").Value = 7512.651899999999
$ws.Range("M122").Value = -5062.651899999999

$ws.Range("H123").Value = 28988.166
$ws.Range("J123").Value = 41571.6
$ws.Range("L123").Value = 41571.6
$ws.Range("N123").Value = -51371.6

$ws.Range("H132").Value = 3490.2632
$ws.Range("I132").Value = 4720.5557
$ws.Range("J132").Value = 2383
$ws.Range("K132").Value = 14161.6671
$ws.Range("L132").Value = 7149
$ws.Range("M132").Value = -11631.6671
$ws.Range("N132").Value = -12209

$ws.Range("H136").Value = 3496.1794
$ws.Range("I136").Value = 3080.6333
$ws.Range("K136").Value = 9241.8999
$ws.Range("M136").Value = -6691.8999
